# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update period labels (column E) for rows 17-26 so they read in ascending
# order (2105, 2106, 2108, 2109, 2110, 2111, 2112, 2201, 2202, 2203) instead
# of the previous descending order.
$ws.Range("E17").Value = "2105"
$ws.Range("E18").Value = "2106"
$ws.Range("E19").Value = "2108"
$ws.Range("E20").Value = "2109"
$ws.Range("E21").Value = "2110"
$ws.Range("E22").Value = "2111"
$ws.Range("E23").Value = "2112"
$ws.Range("E24").Value = "2201"
$ws.Range("E25").Value = "2202"
$ws.Range("E26").Value = "2203"

# Update the "Valor Mora" (F) and "Salario Basico" (G) figures to the new
# totals from the updated account-statement database.
$ws.Range("G16").Value = 908526

$ws.Range("F17").Value = 330703
$ws.Range("G17").Value = 11810838

$ws.Range("G18").Value = 11810838
$ws.Range("G19").Value = 11810838
$ws.Range("G20").Value = 11810838
$ws.Range("G21").Value = 11810838
$ws.Range("G22").Value = 11810838
$ws.Range("G23").Value = 11810838
$ws.Range("G24").Value = 11810838
$ws.Range("G25").Value = 11810838

$ws.Range("F26").Value = 297633
$ws.Range("G26").Value = 11810838
